$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing cell contents (keeps formatting) so the shared-string table
# rebuilds cleanly in the exact order we write it, matching the target layout.
$ws.Cells.ClearContents()

$ws.Cells.Item(1, 1).Value = "Cluster name"
$ws.Cells.Item(1, 2).Value = "Active cases"

$names = @(
    '126 Racecourse Road Public Housing Tower Flemington',
    '3535 Opal Meadow Heights Aged Care Community Meadow Heights',
    '95 Napier Street Apartment Complex Fitzroy',
    'Al Haj Halal Meats Glenroy',
    'Al-Taqwa College Truganina',
    'Amiga Montessori Craigieburn',
    'Apartment Complex 120 Racecourse Road North Melbourne',
    'Australia Post Distribution Centre Sunshine West',
    'Australian Lamb Colac East',
    'Baxter Foods Australia Campbellfield',
    'Budget Car and Truck Rentals Campbellfield',
    'CFMEU Melbourne Office',
    'CS Square Caroline Springs',
    'Cafe Roco Dandenong',
    'Campbellfield Ford Complex Vaccination Clinic Campbellfield',
    'Caroline Springs Police Station',
    'Cedars Medical Clinic Coburg',
    'Chemist Warehouse Campbellfield DC',
    'Chemist Warehouse Fillo Drive Somerton',
    'City of Wyndham Community',
    'Coles Campbellfield Plaza Campbellfield',
    'Coles Coburg North Village',
    'Coles Pakenham Place Shopping Centre',
    'Coles Roxburgh Village Roxburgh Park',
    'Community Kids Bayswater Early Education Centre Bayswater North',
    'Construction Site 1 Warde Street Footscray',
    'Construction Site Olea Apartment Caulfield North',
    'Costco Wholesale Epping',
    'Crusader Caravans Epping',
    'Dandenong Police Station Dandenong',
    'DayHab Rehabilitation Treatment Centre Ringwood East',
    'Direct Freight Express Campbellfield',
    'Disability Residence Life without Barriers Ashwood',
    'Don Watson Coldstore Derrimut',
    'Epworth Healthcare Epworth Richmond Emergency Department',
    'Ermha365 Ltd Doveton',
    'FedEx Station Melbourne Airport',
    'Fine Food Holdings Pty Ltd Dandenong South',
    'Fitzroy Community School Fitzroy North',
    'Fonterra Manufacturing Workplace Campbellfield',
    'General Foods Campbellfield',
    'Gladstone Parade Early Learning & Kinder Glenroy',
    'Goodstart Early Learning Altona',
    'Green Leaves Early Learning Cairnlea',
    'Green Leaves Early Learning Centre Highlands Craigieburn',
    'Greenvale Primary School',
    'Hamilton Marino 236 Jasper Road McKinnon',
    'Hello Fresh Warehouse Ravenhall',
    'Hickory Construction Site Chadstone Car Park Malvern East',
    'IGA Meadow Heights Shopping Centre Meadow Heights',
    'ISS Factory Level 1 Terminal 2 Melbourne Airport Tullamarine',
    'Ibis Kingsgate Hotel Melbourne',
    'Ilim Learning Sanctuary Glenroy',
    'Industrial Galvanizers Valmont Coatings Campbellfield',
    'Inghams Enterprises Thomastown',
    'Kasr Sweets Coolaroo',
    'Kippers Seafood Werribee',
    'Kool Kidz Childcare Narre Warren',
    'Level Crossing Removal Project Lilydale Construction Site John Street',
    'Lineage Logistics Laverton North',
    'Linfox Somerton National Distribution Centre Somerton',
    'McDonald''s Craigieburn North',
    'Mecca D.C Warehouse Melbourne Airport',
    'Melbourne Assessment Prison West Melbourne',
    'Melbourne Metropolitan Remand Centre Ravenhall',
    'Melbourne West Police Station Docklands',
    'Mill Park Police Station Mill Park',
    'MyCentre Childcare Broadmeadows',
    'National Gallery of Victoria Melbourne',
    'Nido Early School Ascot Vale',
    'Nido Early School Glenroy',
    'Northern Health Northern Hospital Epping Emergency Department Tier 1B',
    'Northern Health The Northern Hospital Epping',
    'OnQ Plumbing and Excavations Craigieburn',
    'Oporto Coolaroo',
    'Oscar Romero Catholic Primary School Craigieburn',
    'Our Lady Help of Christian''s Primary School Brunswick East',
    'Pacific Meat Thomastown',
    'Panorama Construction Site Whitehorse Rd Box Hill',
    'Ramsay Health Care Warringal Private Hospital Heidelberg',
    'Ravenhall Correctional Centre Ravenhall',
    'Richmond Quarter 261-271 Bridge Road Construction Site Richmond',
    'Sacca''s Fruit World Broadmeadows Central Shopping Centre',
    'Sharpline Stainless Steel Coburg North',
    'St Margaret''s Primary School OSHC Maribyrnong',
    'St Vincents Hospital Emergency Department Melbourne',
    'Tek Foods Somerton',
    'The Huntly-Goornong Rail Works',
    'The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B',
    'The Royal Melbourne Hospital AMU Ward Parkville',
    'The Royal Talbot Rehabilitation Centre Kew',
    'ThorwestenCabinets Pakenham',
    'Truganina Early Learning Centre Truganina',
    'Unilodge College Square Student Accommodation 570 Lygon Street Carlton',
    'Wallaby Childcare Wollert',
    'Werribee Mercy Hospital Emergency Department',
    'Western Health Footscray Hospital Ward 3B Footscray',
    'Western Health Sunshine Hospital Emergency Department',
    'Yara Childcare Centre Truganina'
)

$values = @(
    5,
    26,
    5,
    14,
    8,
    7,
    5,
    5,
    5,
    5,
    5,
    5,
    11,
    6,
    9,
    9,
    15,
    6,
    11,
    5,
    9,
    21,
    7,
    5,
    15,
    5,
    16,
    29,
    23,
    5,
    6,
    10,
    5,
    5,
    6,
    9,
    14,
    9,
    5,
    9,
    11,
    7,
    11,
    5,
    16,
    5,
    11,
    5,
    5,
    6,
    11,
    5,
    5,
    14,
    6,
    5,
    6,
    12,
    6,
    7,
    10,
    5,
    8,
    7,
    8,
    5,
    7,
    14,
    9,
    11,
    19,
    66,
    16,
    13,
    9,
    5,
    10,
    5,
    5,
    7,
    8,
    12,
    6,
    6,
    12,
    8,
    12,
    6,
    20,
    22,
    10,
    14,
    5,
    7,
    17,
    13,
    10,
    9,
    7
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
